# Update "horarios" (schedules) workbook for Línea 141 with the 01:49:17
# scrape run: refresh the "Última actualización" timestamps, bump the
# "Total filas" counts, and append the newly scraped arrival rows.

$wb = $excel.ActiveWorkbook

# --- Sheet "LP1912": two new rows (9 and 10) ---------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 01:49:17"
$ws1.Range("A3").Value = "Total filas: 5"

$ws1.Range("A9").Value = "01:49:17"
$ws1.Range("B9").Value = "03:03"
$ws1.Range("C9").Value = "215_ALUAR"
$ws1.Range("D9").Value = 74
$ws1.Range("E9").Value = "LP1912"

$ws1.Range("A10").Value = "01:49:17"
$ws1.Range("B10").Value = "03:48"
$ws1.Range("C10").Value = "14_ABASTO"
$ws1.Range("D10").Value = 119
$ws1.Range("E10").Value = "LP1912"

# --- Sheet "LP1912-215": one new row (8) --------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 01:49:17"
$ws2.Range("A3").Value = "Total filas: 3"

$ws2.Range("A8").Value = "01:49:17"
$ws2.Range("B8").Value = "03:03"
$ws2.Range("C8").Value = "215_ALUAR"
$ws2.Range("D8").Value = 74
$ws2.Range("E8").Value = "LP1912"

# --- Sheet "6203-6173": timestamp refresh only --------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 01:49:17"
